$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: was blank in col A / "9146830 - Danúbia Caporusso Bargos" in B/C ("Docentes responsáveis:" data row).
# Becomes the "Programa resumido:" label row, with "Semestral" as its value, height 60.
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 14: "Short syllabus:" label, unchanged text, stays height 60.
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Theory elements and history of urban planning. Theory and practice of environmental planning; environmental planning as an inducer of sustainable development; environmental theory applied to urban problems; legislation and urban environmental policy."
$ws.Range("C14").Value = "Theory elements and history of urban planning. Theory and practice of environmental planning; environmental planning as an inducer of sustainable development; environmental theory applied to urban problems; legislation and urban environmental policy."
$ws.Rows.Item(14).RowHeight = 60

# Row 15: "Programa:" label, value becomes "01/01/2018" (reused string), height grows to 120.
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"
$ws.Rows.Item(15).RowHeight = 120

# Row 16: "Syllabus:" label, unchanged text, stays height 120.
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Environmental planning and management introduction. Planning theory and practice origins. Nature of planning and its relations with geography, politics, economy, society, culture and environment. Analyzes, studies and propositions related to different types of urban growth and expansion; Elements for city environmental structuring; Stages, structures and instruments of environmental planning; Environmental indicators and planning; Public participation in environmental planning; National Policy of the Environment; National System of Conservation Units; City Statute; Environmental Zoning; EIA and EIV as innovative instruments; New concepts and principles of urban-environmental master plans;."
$ws.Range("C16").Value = "Environmental planning and management introduction. Planning theory and practice origins. Nature of planning and its relations with geography, politics, economy, society, culture and environment. Analyzes, studies and propositions related to different types of urban growth and expansion; Elements for city environmental structuring; Stages, structures and instruments of environmental planning; Environmental indicators and planning; Public participation in environmental planning; National Policy of the Environment; National System of Conservation Units; City Statute; Environmental Zoning; EIA and EIV as innovative instruments; New concepts and principles of urban-environmental master plans;."
$ws.Rows.Item(16).RowHeight = 120

# Row 17: "Avaliação:" label only (B/C cleared), default height.
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Rows.Item(17).RowHeight = 15

# Row 18: "Método:" label, value becomes "9146830 - Danúbia Caporusso Bargos" (reused string), height 60.
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Range("C18").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Rows.Item(18).RowHeight = 60

# Row 19: "Critério:" label, value is the "Aulas teóricas..." text, height 60.
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas teóricas e práticas, visitas técnicas e exercícios dirigidos. Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios."
$ws.Range("C19").Value = "Aulas teóricas e práticas, visitas técnicas e exercícios dirigidos. Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios."
$ws.Rows.Item(19).RowHeight = 60

# Row 20: "Norma de recuperação:" label, value is the "Média ponderada..." text, height 60.
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios."
$ws.Range("C20").Value = "Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios."
$ws.Rows.Item(20).RowHeight = 60

# Row 21: "Bibliografia:" label, value is the "Provas e/ou exercícios dirigidos." text, height grows to 120.
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Provas e/ou exercícios dirigidos."
$ws.Range("C21").Value = "Provas e/ou exercícios dirigidos."
$ws.Rows.Item(21).RowHeight = 120

# Row 22: "Requisitos:" label only (B/C cleared), default height.
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Rows.Item(22).RowHeight = 15

# Row 23: requirement text moves up from row 24 (label column A cleared), height 30.
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = $ws.Range("B24").Value
$ws.Range("C23").Value = $ws.Range("C24").Value
$ws.Rows.Item(23).RowHeight = 30

# Row 24 is no longer used; clear it entirely and delete it so the used range shrinks to A1:C23.
$ws.Rows.Item(24).Delete()
